# The "Air Cooler" variety prompt (rows 46-58, column B) is being updated
# to mention the newly added "Window Cooler" option.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers 4. Window Cooler"

$target = $ws.Range("B46:B58")
$target.Value = $newText

$target.Select()
